$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.765933666666667
$ws.Range("H2").Value = 5.297801000000001
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.02564166666666666
$ws.Range("N2").Value = 0.076925
$ws.Range("O2").Value = 0.0006780701807970013
$ws.Range("P2").Value = 0.0006780701807970013
$ws.Range("Q2").Value = 0.04528148243611111
$ws.Range("R2").Value = 0.407533341925
$ws.Range("S2").Value = 0.0006780701807970013
$ws.Range("T2").Value = 0.0006780701807970013

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.765933666666667
$ws.Range("H3").Value = 5.297801000000001
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01112833333333333
$ws.Range("N3").Value = 0.033385
$ws.Range("O3").Value = 0.0002942784918545062
$ws.Range("P3").Value = 0.0002942784918545062
$ws.Range("Q3").Value = 0.01965189848722222
$ws.Range("R3").Value = 0.176867086385
$ws.Range("S3").Value = 0.0002942784918545062
$ws.Range("T3").Value = 0.0002942784918545062

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.765933666666667
$ws.Range("H4").Value = 5.297801000000001
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03487066666666667
$ws.Range("N4").Value = 0.104612
$ws.Range("O4").Value = 0.0009221225577320236
$ws.Range("P4").Value = 0.0009221225577320235
$ws.Range("Q4").Value = 0.06157928424577779
$ws.Range("R4").Value = 0.554213558212
$ws.Range("S4").Value = 0.0009221225577320236
$ws.Range("T4").Value = 0.0009221225577320235

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.765933666666667
$ws.Range("H5").Value = 5.297801000000001
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 37.74401233333333
$ws.Range("N5").Value = 113.232037
$ws.Range("O5").Value = 0.9981055287696164
$ws.Range("P5").Value = 0.9981055287696164
$ws.Range("Q5").Value = 66.65342209451522
$ws.Range("R5").Value = 599.880798850637
$ws.Range("S5").Value = 0.9981055287696164
$ws.Range("T5").Value = 0.9981055287696164

